# estop-1.0-bom.xlsx: "changed antenna to right angle and removed crystal"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 was "220uF electrolytic" / "732-8911-1-ND" -> now "47uF ceramic" / "587-1780-1-ND",
# with a Unit Cost of 0.7 added.
$ws.Range("A6").Value = "47uF ceramic"
$ws.Range("B6").Value = "587-1780-1-ND"
$ws.Range("E6").Value = 0.7

# Row 7's "Order" quantity (D7 = 2) is cleared.
$ws.Range("D7").ClearContents()

# The "8MHz crystal" (row 12) and "18pF ceramic cap" (row 13) rows are removed entirely,
# shifting the Level shifter / SMA connector rows up to 12 / 13.
$ws.Rows("12:13").Delete() | Out-Null

# Update the saved selection to match the post-edit cursor position.
$ws.Range("A12:XFD12").Select() | Out-Null
